$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '93.701.84'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +3.76%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.123.27'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -1.38%  '

$ws.Range("E4").Value = '  -0.06%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '243.54'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +2.41%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '617.15'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.84%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '1.09'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -1.02%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.414'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +11.94%  '

$ws.Range("E9").Value = '  -0.14%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '3.122.45'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -1.30%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.737'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -0.34%  '

$ws.Range("E12").Value = '  -0.36%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.0000258'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +4.55%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '34.77'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -1.79%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '92.890.87'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +2.67%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '5.49'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.62%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '3.711.72'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -1.41%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.120.91'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -1.96%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '3.80'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +3.07%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '14.81'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -1.59%  '

$ws.Range("E21").Value = '  +2.81%  '

$ws.Range("E22").Value = '  -0.70%  '

$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '9.45'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +3.72%  '

$ws.Range("B24").Value = 'BitcoinCash'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '452.06'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +2.57%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '5.84'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +1.67%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '87.58'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -1.64%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '11.88'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -1.09%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '3.292.11'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.13%  '

$ws.Range("E29").Value = '  -0.06%  '

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.136'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +6.55%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '0.169'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -0.17%  '

$ws.Range("E32").Value = '  -1.40%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '9.25'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -2.60%  '

$ws.Range("E34").Value = '  +0.29%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '8.13'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +4.49%  '

$ws.Range("E36").Value = '  -2.47%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '26.27'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.22%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '3.98'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +6.01%  '

$ws.Range("E39").Value = '  -1.80%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '483.25'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -4.34%  '

$ws.Range("E41").Value = '  -3.43%  '

$ws.Range("E42").Value = '  +3.24%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.438'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.75%  '

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '23.12'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +4.63%  '

$ws.Range("E45").Value = '  +0.02%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '162.39'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +3.34%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.94'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +1.30%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '0.695'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -3.65%  '

$ws.Range("E49").Value = '  +2.00%  '

$ws.Range("E50").Value = '  +3.57%  '

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '4.47'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +1.11%  '
